$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '21.744.65'
$ws.Range("E2").Value = '  -1.16%  '
$ws.Range("D3").Value = '1.541.23'
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9981'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9993'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '290.37'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.89%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3949'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.82%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3208'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.18'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07217'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.087'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9983'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.783'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.52'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.663'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.43%  '
$ws.Range("D16").Value = '1.542.80'
$ws.Range("E16").Value = '  -0.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001104'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06614'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '84.43'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9988'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.173'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.59%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.66'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.92'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.370'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.53%  '
$ws.Range("D25").Value = '21.744.50'
$ws.Range("E25").Value = '  -1.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.420'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.49%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '152.32'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.24%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.58'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.876'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.86%  '
$ws.Range("D30").Value = '1.714.46'
$ws.Range("E30").Value = '  -0.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '117.90'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.63%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.168'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +8.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9815'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08159'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.651'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.245'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02256'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06047'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.55%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.493'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.43'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2057'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.188'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9986'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5869'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.23'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.58%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.737'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5620'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.33%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.910'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.18%  '
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.174'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.50%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '117.30'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06755'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.94%  '